# Actualización automática 2025-07-10 08:45:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("E11").Value = 565.27
$wsGrupo.Range("E12").Value = "1 de 10"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F11").Value = 565.27
$wsMensual.Range("F12").Value = 705.64
# Column F width: 11 -> 12 (character units require a slight offset to land on 12)
$wsMensual.Columns.Item(6).ColumnWidth = 11.14

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D2").Value = 2055.94
$wsCumpl.Range("E2").Value = -2055.94
$wsCumpl.Range("D4").Value = 2967.62
$wsCumpl.Range("E4").Value = 10755.72
$wsCumpl.Range("F4").Value = 0.2162461907961181
